$wb = $excel.ActiveWorkbook

# --- Instructions sheet: update two help-text cells ---------------------
$instructions = $wb.Worksheets.Item("Instructions")

$instructions.Range("A12").Value = "It is highly recommended that you ensure your website is backed up before importing from your CSV files."

$instructions.Range("A41").Value = "The ‘Category’ must be an existing event category.  Create new event categories if necessary before importing events that refer to them.  Multiple event categories should be separated using the vertical bar character |,  for example ‘Meeting|Social’."

# --- Instructions sheet: move the saved view/selection back to the top --
$instructions.Activate() | Out-Null
$instructions.Range("A12").Select() | Out-Null
